# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Recomputed strikeout (K) values for column G replace the old Strike# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 0
    6  = 2
    7  = 1
    8  = 1
    9  = 1
    10 = 1
    11 = 0
    12 = 0
    13 = 1
    14 = 0
    15 = 1
    16 = 1
    17 = 0
    18 = 2
    19 = 2
    20 = 0
    21 = 1
    22 = 1
    23 = 2
    24 = 1
    25 = 1
    26 = 1
    27 = 0
    28 = 1
    29 = 2
    30 = 0
    31 = 1
    32 = 0
    33 = 1
    34 = 1
    35 = 1
    36 = 0
    37 = 1
    38 = 2
    39 = 1
    40 = 0
    41 = 1
    42 = 1
    43 = 0
    44 = 2
    45 = 3
    46 = 1
    47 = 2
    48 = 1
    49 = 1
    50 = 0
    51 = 2
    52 = 1
    53 = 1
    54 = 0
    55 = 2
    56 = 0
    57 = 1
    58 = 1
    62 = 2
    64 = 1
    65 = 1
    66 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
